# Auto update Excel log - append new mmWave sensor readings

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

$newRows = @(
    @("2026-01-30", "15:06:13", "15:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-30", "15:06:26", "15:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-30", "15:06:37", "15:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-30", "15:06:47", "15:00", "Living Room", "PRESENCE_DETECTED", "Active")
)

$startRow = 54
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]

    # Column A holds a date-like string ("2026-01-30"). Assigning it directly
    # via .Value makes Excel auto-convert it into a date serial number, so we
    # temporarily force text format, assign the value, then restore the
    # cell's style back to Normal so no stray formatting is left behind.
    $dateCell = $ws.Cells.Item($row, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $data[0]
    $dateCell.Style = "Normal"

    for ($col = 2; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $data[$col - 1]
    }
}
